# refatoração - cálculos de apoio médio
#
# Splits the former "*_sucesso" stats for arrecadado into explicit
# arrecadado_avg/std/min/max columns, adds std/min/max alongside the
# existing apoio_medio column, and adds std/min/max alongside the
# renamed contribuicoes_med column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert 3 new columns right after "apoio_medio" (L) for
#    apoio_std / apoio_min / apoio_max. They inherit the R$ number
#    format from column L, same as the arrecadado_* block.
$ws.Columns("M:O").Insert()

# 2) Insert 3 new columns right after "contribuicoes_med" (now Q) for
#    contribuicoes_std / contribuicoes_min / contribuicoes_max. They
#    inherit the #,##0 number format from column Q.
$ws.Columns("R:T").Insert()

# Rename the old *_sucesso headers to arrecadado_* (same columns/position)
$ws.Range("H1").Value = "arrecadado_avg"
$ws.Range("I1").Value = "arrecadado_std"
$ws.Range("J1").Value = "arrecadado_min"
$ws.Range("K1").Value = "arrecadado_max"

# New apoio_* headers in the freshly inserted columns
$ws.Range("M1").Value = "apoio_std"
$ws.Range("N1").Value = "apoio_min"
$ws.Range("O1").Value = "apoio_max"

# Rename media_contribuicoes -> contribuicoes_med (shifted to Q)
$ws.Range("Q1").Value = "contribuicoes_med"

# New contribuicoes_* headers in the freshly inserted columns
$ws.Range("R1").Value = "contribuicoes_std"
$ws.Range("S1").Value = "contribuicoes_min"
$ws.Range("T1").Value = "contribuicoes_max"

# Updated apoio_medio value (recomputed)
$ws.Range("L2").Value = 91.85574933975617

# New apoio_std / apoio_min / apoio_max values
$ws.Range("M2").Value = 49.08980856017526
$ws.Range("N2").Value = 13.93896149503088
$ws.Range("O2").Value = 792.0360759681182

# New contribuicoes_std / contribuicoes_min / contribuicoes_max values
$ws.Range("R2").Value = 423.019225146675
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 6494
